$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Averaged")

# Update relative_age (D) and absolute_age (E) values for rows 3-36
# to reflect the rerun calculations (commit: "Reran calculations").
$ws.Range("D3").Value = 0.035223931546153893
$ws.Range("E3").Value = 201.52877606845385
$ws.Range("D4").Value = 0.050868432627401879
$ws.Range("E4").Value = 201.51313156737257
$ws.Range("D5").Value = 0.073624070563762758
$ws.Range("E5").Value = 201.49037592943623
$ws.Range("D6").Value = 0.10320639988103184
$ws.Range("E6").Value = 201.46079360011896
$ws.Range("D7").Value = 0.10605085462307688
$ws.Range("E7").Value = 201.45794914537692
$ws.Range("D8").Value = 0.1074730819940995
$ws.Range("E8").Value = 201.45652691800589
$ws.Range("D9").Value = 0.11117087315875811
$ws.Range("E9").Value = 201.45282912684124
$ws.Range("D10").Value = 0.11919290281584802
$ws.Range("E10").Value = 201.44480709718414
$ws.Range("D11").Value = 0.12458510624494902
$ws.Range("E11").Value = 201.43941489375504
$ws.Range("D12").Value = 0.13195445093138711
$ws.Range("E12").Value = 201.43204554906859
$ws.Range("D13").Value = 0.13590873344606114
$ws.Range("E13").Value = 201.42809126655393
$ws.Range("D14").Value = 0.14004275607503866
$ws.Range("E14").Value = 201.42395724392495
$ws.Range("D15").Value = 0.14058197641794873
$ws.Range("E15").Value = 201.42341802358203
$ws.Range("D16").Value = 0.14525521938983624
$ws.Range("E16").Value = 201.41874478061015
$ws.Range("D17").Value = 0.14885002167590358
$ws.Range("E17").Value = 201.41514997832408
$ws.Range("D18").Value = 0.14938924201881371
$ws.Range("E18").Value = 201.41461075798117
$ws.Range("D19").Value = 0.15208534373336416
$ws.Range("E19").Value = 201.41191465626662
$ws.Range("D20").Value = 0.15406248499070119
$ws.Range("E20").Value = 201.40993751500929
$ws.Range("D21").Value = 0.1560396262480383
$ws.Range("E21").Value = 201.40796037375196
$ws.Range("D22").Value = 0.16826195402066721
$ws.Range("E22").Value = 201.39573804597933
$ws.Range("D23").Value = 0.17167701619243128
$ws.Range("E23").Value = 201.39232298380756
$ws.Range("D24").Value = 0.17599077893571205
$ws.Range("E24").Value = 201.38800922106429
$ws.Range("D25").Value = 0.21234619175508065
$ws.Range("E25").Value = 201.35165380824492
$ws.Range("D26").Value = 0.21321871787094013
$ws.Range("E26").Value = 201.35078128212905
$ws.Range("D27").Value = 0.21601833800389364
$ws.Range("E27").Value = 201.3479816619961
$ws.Range("D28").Value = 0.22826667608556503
$ws.Range("E28").Value = 201.33573332391444
$ws.Range("D29").Value = 0.25171349469905036
$ws.Range("E29").Value = 201.31228650530093
$ws.Range("D30").Value = 0.27568524208746442
$ws.Range("E30").Value = 201.28831475791253
$ws.Range("D31").Value = 0.27743500467056026
$ws.Range("E31").Value = 201.28656499532943
$ws.Range("D32").Value = 0.31715154688078667
$ws.Range("E32").Value = 201.24684845311921
$ws.Range("D33").Value = 0.3317449808455657
$ws.Range("E33").Value = 201.23225501915442
$ws.Range("D34").Value = 0.36056239475069901
$ws.Range("E34").Value = 201.20343760524929
$ws.Range("D35").Value = 0.36795147523919464
$ws.Range("E35").Value = 201.19604852476081
$ws.Range("D36").Value = 0.37404746664220362
$ws.Range("E36").Value = 201.1899525333578

# Widen column D slightly to match the updated layout used for the evolutions figure.
$ws.Range("D1").EntireColumn.ColumnWidth = 13.7109375
